$wb = $excel.ActiveWorkbook

# --- Sheet "Python": update C29 text and remove trailing rows 30-37 ---
$wsPython = $wb.Worksheets.Item("Python")
$wsPython.Range("C29").Value = "Загрузка CSV-данных в датафрейм"
$wsPython.Rows.Item(30).Resize(8).EntireRow.Delete()

# --- Sheet "Links": drop the old "Яндекс" row (row 2) and replace the ---
# --- former top row (Google search engine) with a new CLI link entry ---
$wsLinks = $wb.Worksheets.Item("Links")
$wsLinks.Rows.Item(2).EntireRow.Delete()
$wsLinks.Range("A1").Value = 33
$wsLinks.Range("B1").Value = "Что такое CLI?"
$wsLinks.Range("C1").Value = "ru.wikipedia.org/wiki/Интерфейс_командной_строки"
$wsLinks.Range("D1").Value = "Что такое CLI?"

# --- Sheet "Bash": fix wording in C74 ---
$wsBash = $wb.Worksheets.Item("Bash")
$wsBash.Range("C74").Value = "Перенос файла на другой диск"
